# Apply updated crypto price / volume(1h) figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'72.103.34"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.96%  '

$ws.Range('D3').Value = "'2.665.44"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.39%  '

$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').Value = "'597.85"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.38%  '

$ws.Range('D6').Value = "'175.66"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.66%  '

$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').Value = "'0.522"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.68%  '

$ws.Range('D9').Value = "'2.662.81"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.26%  '

$ws.Range('E10').Value = '  +2.15%  '

$ws.Range('E11').Value = '  +2.28%  '

$ws.Range('E12').Value = '  +1.54%  '

$ws.Range('D13').Value = "'4.98"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.49%  '

$ws.Range('D14').Value = "'3.149.74"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.23%  '

$ws.Range('E15').Value = '  -1.16%  '

$ws.Range('D16').Value = "'71.953.23"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.86%  '

$ws.Range('D17').Value = "'26.18"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.22%  '

$ws.Range('D18').Value = "'2.704.15"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.28%  '

$ws.Range('D19').Value = "'12.03"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.68%  '

$ws.Range('D20').Value = "'7.96"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.07%  '

$ws.Range('D21').Value = "'370.40"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.19%  '

$ws.Range('E22').Value = '  -0.29%  '

$ws.Range('D23').Value = "'2.02"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.03%  '

$ws.Range('D24').Value = "'71.70"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.98%  '

$ws.Range('E25').Value = '  +0.04%  '

$ws.Range('D26').Value = "'4.30"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.74%  '

$ws.Range('D27').Value = "'9.79"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.88%  '

$ws.Range('D28').Value = "'2.799.05"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.15%  '

$ws.Range('D29').Value = "'0.999"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.12%  '

$ws.Range('D30').Value = "'0.0₃0936"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.79%  '

$ws.Range('E31').Value = '  -0.82%  '

$ws.Range('D32').Value = "'506.82"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.69%  '

$ws.Range('D33').Value = "'1.29"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.71%  '

$ws.Range('E34').Value = '  -1.27%  '

$ws.Range('D35').Value = "'0.999"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.16%  '

$ws.Range('D36').Value = "'163.60"
$ws.Range('D36').Style = 'Normal'

$ws.Range('D37').Value = "'19.46"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.94%  '

$ws.Range('D38').Value = "'19.08"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.20%  '

$ws.Range('E39').Value = '  -1.28%  '

$ws.Range('E40').Value = '  -5.09%  '

$ws.Range('E41').Value = '  -8.37%  '

$ws.Range('E42').Value = '  -0.04%  '

$ws.Range('E43').Value = '  -1.36%  '

$ws.Range('D44').Value = "'2.54"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.15%  '

$ws.Range('E45').Value = '  -0.24%  '

$ws.Range('D46').Value = "'39.17"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.93%  '

$ws.Range('D47').Value = "'153.05"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.66%  '

$ws.Range('D48').Value = "'3.70"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.47%  '

$ws.Range('D49').Value = "'0.546"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.68%  '

$ws.Range('E50').Value = '  +1.86%  '

$ws.Range('D51').Value = "'0.0759"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.48%  '
